# Apply "repull data" updates to column F (dSF) for several rows.
# This mirrors re-pulled values overwriting the previously stored dSF figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    7  = 5
    8  = -3
    15 = -2
    23 = -3
    24 = -2
    27 = -2
    31 = -1
    32 = -3
    38 = 2
    41 = -4
    45 = -8
    46 = 5
    47 = 0
    49 = -2
    50 = 0
    53 = 3
    58 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
